{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst TXT_ROUTE_CREATE = \"Route Create \u2013 ostatni magazyn, by by\u0142 rzeczywi\u015bcie ostatni (ale jeszcze do przemy\u015blenia), mo\u017ce po prostu do usuni\u0119cia\";\nconst TXT_ROUTE_DETAILS_CURVE = \"Route Details \u2013 poprawa algorytmu wyznaczaj\u0105cego trasy, tak by bra\u0142 pod uwag\u0119 krzywizn\u0119 ziemi \u2013 algorytmy w internecie\";\nconst TXT_RESOURCES_CREATE_PREFIX = \"Resources Create \u2013 dodanie tabeli\";\n\nlet routeCreatePara = null;\nlet routeDetailsCurvePara = null;\nlet resourcesCreatePara = null;\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const p = paras.items[i];\n  const t = p.text;\n  if (t === TXT_ROUTE_CREATE) {\n    routeCreatePara = p;\n  } else if (t === TXT_ROUTE_DETAILS_CURVE) {\n    routeDetailsCurvePara = p;\n  } else if (t.indexOf(TXT_RESOURCES_CREATE_PREFIX) === 0) {\n    resourcesCreatePara = p;\n  }\n}\n\n// 1) Delete the \"Route Create\" paragraph entirely.\nif (routeCreatePara) {\n  routeCreatePara.delete();\n}\n\n// 2) Delete the \"Route Details - krzywizna ziemi\" paragraph entirely.\nif (routeDetailsCurvePara) {\n  routeDetailsCurvePara.delete();\n}\n\nawait context.sync();\n\n// 3) Remove the old \"_GoBack\" bookmark that used to sit at the end of the\n//    \"Drugim zadaniem\" paragraph; it is about to be re-created around\n//    \"dodawania\" inside the \"Resources Create\" paragraph below (bookmark\n//    names must stay unique in the document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) Split the \"Resources Create\" paragraph's run so the word \"dodawania\"\n//    sits inside its own bookmark named \"_GoBack\".\nif (resourcesCreatePara) {\n  const theWord = \"dodawania\";\n\n  const hits = resourcesCreatePara.search(theWord, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items/text\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    const wordRange = hits.items[0];\n    wordRange.insertBookmark(\"_GoBack\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$routeCreateText = \"Route Create \u2013 ostatni magazyn, by by\u0142 rzeczywi\u015bcie ostatni (ale jeszcze do przemy\u015blenia), mo\u017ce po prostu do usuni\u0119cia\"\n$routeDetailsCurveText = \"Route Details \u2013 poprawa algorytmu wyznaczaj\u0105cego trasy, tak by bra\u0142 pod uwag\u0119 krzywizn\u0119 ziemi \u2013 algorytmy w internecie\"\n$resourcesCreatePrefix = \"Resources Create \u2013 dodanie tabeli\"\n\n# 1) Delete the \"Route Create\" paragraph entirely.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq $routeCreateText) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Delete the \"Route Details - krzywizna ziemi\" paragraph entirely.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq $routeDetailsCurveText) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 3) Remove the old \"_GoBack\" bookmark sitting at the end of the \"Drugim\n#    zadaniem\" paragraph; it is about to be re-created around \"dodawania\"\n#    inside the \"Resources Create\" paragraph (bookmark names must stay\n#    unique in the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 4) Wrap the word \"dodawania\" (inside \"dodawania/aktualizowania\") in the\n#    \"Resources Create\" paragraph with a new \"_GoBack\" bookmark.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith($resourcesCreatePrefix)) {\n        $searchRange = $p.Range.Duplicate\n        $find = $searchRange.Find\n        $find.ClearFormatting()\n        $find.Text = \"dodawania\"\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.Forward = $true\n        $find.Wrap = 0\n        if ($find.Execute()) {\n            $d.Bookmarks.Add(\"_GoBack\", $searchRange)\n        }\n        break\n    }\n}\n"}
